$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 245; this shifts existing rows 245-320 down
# to 246-321 (matching the dimension change from A1:R320 to A1:R321) and
# keeps every other row's data untouched.
$ws.Rows(245).Insert()

# Populate the newly inserted row 245 with the new weekly price entry.
$ws.Cells.Item(245, 1).Value  = 4
$ws.Cells.Item(245, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(245, 3).Value  = "Los Lagos"
$ws.Cells.Item(245, 4).Value  = 44588
$ws.Cells.Item(245, 5).Value  = 10
$ws.Cells.Item(245, 6).Value  = 100114001
$ws.Cells.Item(245, 7).Value  = "Papa"
$ws.Cells.Item(245, 8).Value  = "Patagonia"
$ws.Cells.Item(245, 9).Value  = "1a nueva(o)"
$ws.Cells.Item(245, 10).Value = 250
$ws.Cells.Item(245, 11).Value = 9000
$ws.Cells.Item(245, 12).Value = 10000
$ws.Cells.Item(245, 13).Value = 9600
$ws.Cells.Item(245, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(245, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(245, 16).Value = 384
$ws.Cells.Item(245, 17).Value = 25
$ws.Cells.Item(245, 18).Value = "Hortaliza"
